$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-08 Thursday" "2024-08-09 Friday"

Replace-Text "940×8=" "694×3="
Replace-Text "399×2=" "102×4="
Replace-Text "958×2=" "577×8="
Replace-Text "770×3=" "446×4="
Replace-Text "512×6=" "276×5="
Replace-Text "686×4=" "528×9="
Replace-Text "454×9=" "435×2="
Replace-Text "573×9=" "799×2="
Replace-Text "649×3=" "599×2="
Replace-Text "613×9=" "691×3="
Replace-Text "814×9=" "211×2="
Replace-Text "520×4=" "380×9="
Replace-Text "926×2=" "693×3="
Replace-Text "633×3=" "115×2="
Replace-Text "713×8=" "350×4="
Replace-Text "502×9=" "373×2="
Replace-Text "120×6=" "190×3="
Replace-Text "452×5=" "544×9="
Replace-Text "133×8=" "673×6="
Replace-Text "239×3=" "885×5="
Replace-Text "648×2=" "869×2="
Replace-Text "571×7=" "463×3="
Replace-Text "477×9=" "803×5="
Replace-Text "769×5=" "793×8="
Replace-Text "145×3=" "776×6="
